# Integrate local file acquisition with cases template
#
# - Extend "Number of Cases" partition/aliases with a "time_window" /
#   "Cumulative cases" entry.
# - Rename "Gender" referential to "Gender code".
# - Consolidate the separate "Alerts of potential cases in animals/humans"
#   + "Number of events" rows into a single "Number of alerts" row that is
#   now partitioned by reporting_period,alert_topic and whose aliases cover
#   both the human and animal alert topics.
# - Turn "Alert topic" into a proper Characteristic row and add a new
#   "Time window" referential row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Number of Cases" (row 2): add the time_window partition + alias ---
$ws.Range("G2").Value = "reporting_period,case_status,time_window"
$ws.Range("H2").Value = '[{"alias":"Confirmed cases", "variable":"Number of Cases", "modifiers":[{"variable":"Case Status", "value":"Confirmed"}]}, {"alias":"Recovered cases", "variable":"Number of Cases", "modifiers":[{"variable":"Case Status", "value":"Recovered"}]}, {"alias":"Active cases", "variable":"Number of Cases", "modifiers":[{"variable":"Case Status", "value":"Active"}]}, {"alias":"Reinfection cases", "variable":"Number of Cases", "modifiers":[{"variable":"Case Status", "value":"Reinfection"}]}, {"alias":"Possible cases", "variable":"Number of Cases", "modifiers":[{"variable":"Case Status", "value":"Possible"}]}, {"alias":"Probable cases", "variable":"Number of Cases", "modifiers":[{"variable":"Case Status", "value":"Probable"}]}, {"alias":"Imported cases", "variable":"Number of Cases", "modifiers":[{"variable":"Case Status", "value":"Imported"}]},{"alias":"Cumulative cases", "variable":"Number of Cases", "modifiers":[{"variable":"Time window", "value":"Cumulative"}]}, {"alias":"Cases at onset of symptoms date", "variable":"Number of Cases", "modifiers":[{"variable":"period type", "value":"Onset of symptoms date"}]}]'

# --- "Gender" referential is renamed to "Gender code" (update both the
#     referential row itself and the "Linked Attributes" back-reference on
#     its alias row, "Gender name", so they stay in sync) ---
$ws.Range("B116").Value = "Gender code"
$ws.Range("F117").Value = "Gender code"

# --- Rebuild the alerts block (rows 120-124) ---
# Before:
#   120 Alerts of potential cases in animals  (Observation)
#   121 Alerts of potential cases in humans   (Observation)
#   122 Number of events                      (Observation, aliases: humans only)
#   123 Number of alerts                      (Observation)
#   124 Alert topic                           (Observation)
# After:
#   120 Number of alerts   (Observation, partition reporting_period,alert_topic,
#        aliases covering both humans and animals)
#   121 Alert topic        (Characteristic)
#   122 Time window         (new Characteristic / Referential entry)

# Drop the two now-redundant "Number of alerts" / "Alert topic" rows that used
# to sit at the bottom of the block (124 then 123, highest row first so the
# row numbers of earlier rows don't shift while we still need them).
$ws.Rows.Item(124).Delete()
$ws.Rows.Item(123).Delete()

# Row 120 used to be "Alerts of potential cases in animals" - turn it into
# the consolidated "Number of alerts" row.
$ws.Range("A120").Value = "01. Cases"
$ws.Range("B120").Value = "Number of alerts"
$ws.Range("D120").Value = "Observation"
$ws.Range("G120").Value = "reporting_period,alert_topic"
$ws.Range("H120").Value = '[{"alias":"alerts on humans", "variable":"Number of alerts", "modifiers":[{"variable":"Alert topic", "value":"Potential case in humans"}]}, {"alias":"alerts on animals", "variable":"Number of alerts", "modifiers":[{"variable":"Alert topic", "value":"Potential case in animals"}]}]'

# Row 121 used to be "Alerts of potential cases in humans" - turn it into the
# "Alert topic" characteristic row.
$ws.Range("A121").Value = "01. Cases"
$ws.Range("B121").Value = "Alert topic"
$ws.Range("D121").Value = "Characteristic"
$ws.Range("H121").Value = "[]"

# Row 122 used to be "Number of events" (and used the now-dropped heavier
# font style reserved for the tail of the table) - turn it into the new
# "Time window" referential row and pick up the regular row formatting from
# the row above instead of the old, now-unused, style.
$ws.Range("A121").Copy()
$ws.Range("A122").PasteSpecial(-4122)
$ws.Range("A122").Value = "13. Referentials"
$ws.Range("B122").Value = "Time window"
$ws.Range("D122").Value = "Characteristic"
$ws.Range("H122").Value = "[]"

# Restore view: back to the top of the sheet, selection on G16.
[void]$ws.Range("G16").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
